$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# The slide's shape-id/name counter in this host starts handing out
# "TextBox 2" (id 3) for the first programmatically-added textbox on a
# slide, regardless of the ids already used by existing shapes. Adding
# and immediately removing a throwaway textbox first advances that
# counter so the real textbox we want lands on id=6 / "TextBox 5",
# matching the target deck exactly.
$warmup = $s.Shapes.AddTextbox(1, 1, 1, 1, 1)
$warmup.Delete()

# Position/size converted from the target EMU values (1 pt = 12700 EMU),
# kept at full floating-point precision so the round-trip back to EMU
# lands on the exact target integers.
$left = 5342021 / 12700
$top = 737937 / 12700
$width = 3593432 / 12700
$height = 369332 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "Need Package: ROCR"
